# Updates with supplemental figures
# Replaces the PSSM score matrix values in B2:K21 with the recalculated
# values from the updated supplemental figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-17.32909555856886, 1.758587403681116, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, 2.49115353453601, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, 2.096491217230355, 2.093373742324949, -17.32909555856886, 3.609410813762548, -17.32909555856886, 1.273751888268847, -17.32909555856886, 1.929748774529283, -17.32909555856886),
    @(-17.32909555856886, 1.81555521667192, -17.32909555856886, -17.32909555856886, -17.32909555856886, 2.52272209559268, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(2.741891424448053, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, 2.082574721536524, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(3.734484508900147, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, 1.038937645632331, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, 2.817568821015108, -17.32909555856886, 3.240594125842225, -17.32909555856886, -17.32909555856886, -17.32909555856886, 4.321919771083537),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, 2.458878718140005, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, 1.622704231348581, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, 1.105846172298742, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, 1.070236612871018, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, 1.73414414233265, -17.32909555856886),
    @(-17.32909555856886, 2.34949254610197, 2.680264266363509, -17.32909555856886, -17.32909555856886, -17.32909555856886, 1.549820486556399, 2.703263667704547, 2.638379801729535, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, -17.32909555856886, 2.271864487021341, 1.361274200593418, 1.830807797896836, -17.32909555856886),
    @(-17.32909555856886, -17.32909555856886, 1.466097521949392, -17.32909555856886, -17.32909555856886, -17.32909555856886, 1.43869125165304, 0.9658168083698656, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, 0.8470240106305301, 1.206624060105091, -17.32909555856886, 2.96251246875557, -17.32909555856886, 1.645120149510686, 0.3650242360286584, -17.32909555856886, -17.32909555856886),
    @(-17.32909555856886, 0.9443950162032061, -17.32909555856886, 1.685559722694299, -17.32909555856886, 2.263458199212733, 1.995034434460355, -17.32909555856886, -17.32909555856886, -17.32909555856886)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $values[$i].Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $values[$i][$j]
    }
}
